# Daily attendance processing - 2026-01-28 08:05:29
# Reorders the "Recorded By" (column G) comma-separated list so that the
# literal token "System" (exact case) is moved to the front of the list,
# while preserving the relative order of the remaining tokens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        $hasSystem = $false
        foreach ($p in $trimmed) {
            if ($p.Equals("System")) {
                $hasSystem = $true
            }
        }

        if ($hasSystem) {
            $rest = @()
            foreach ($p in $trimmed) {
                if (-not $p.Equals("System")) {
                    $rest += $p
                }
            }
            $newParts = @("System") + $rest
            $newVal = [string]::Join(", ", $newParts)

            if (-not $newVal.Equals($val)) {
                $cell.Value = $newVal
            }
        }
    }
}
